$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B15: updated service name
$b15Text = @"
OnlineSavingAccountTypes?SEARCH_TYPE=ALL_ONLINE_SAVING service
"@
$ws.Range("B15").Value = $b15Text

# F15: updated service URL
$f15Text = @"
https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes?SEARCH_TYPE=ALL_ONLINE_SAVING
"@
$ws.Range("F15").Value = $f15Text

# G15: updated sample JSON response
$g15Text = @"
{
    "d": {
        "results": [
            {
                "__metadata": {
                    "id": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(1)",
                    "uri": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(1)",
                    "type": "com.sap.banking.custom.onlinesaving.endpoint.v1_0.beans.OnlineSavingAccountType"
                },
                "Account_Type_ID": 1,
                "Name_EN": "Tiết kiệm thường trực tuyến",
                "Name_VN": "Tiết kiệm thường trực tuyến",
                "Periodic": null,
                "Tandc_Link": null,
                "Status": "0",
                "Introduce_product_name_en": "Enjoy interest based on the principal amount and flexible term",
                "Introduce_product_name_vn": "Hưởng lãi suất theo số tiền gửi và kỳ hạn gửi đa dạng",
                "Link_detail_product": "https://www.vpbank.com.vn/san-pham/gui-tiet-kiem/tiet-kiem-truc-tuyen"
            },
            {
                "__metadata": {
                    "id": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(2)",
                    "uri": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(2)",
                    "type": "com.sap.banking.custom.onlinesaving.endpoint.v1_0.beans.OnlineSavingAccountType"
                },
                "Account_Type_ID": 2,
                "Name_EN": "Bảo toàn thịnh vượng trực tuyến",
                "Name_VN": "Bảo toàn thịnh vượng trực tuyến",
                "Periodic": null,
                "Tandc_Link": null,
                "Status": "0",
                "Introduce_product_name_en": "Preserve interest rates before the fluctuations in interest rates on the market.",
                "Introduce_product_name_vn": "Bảo toàn lãi suất trước những biến động lãi suất trên thị trường",
                "Link_detail_product": "https://www.vpbank.com.vn/san-pham/gui-tiet-kiem/tiet-kiem-bao-toan-thinh-vuong"
            },
            {
                "__metadata": {
                    "id": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(3)",
                    "uri": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(3)",
                    "type": "com.sap.banking.custom.onlinesaving.endpoint.v1_0.beans.OnlineSavingAccountType"
                },
                "Account_Type_ID": 3,
                "Name_EN": "Tiết kiệm lĩnh lãi định kỳ trực tuyến",
                "Name_VN": "Tiết kiệm lĩnh lãi định kỳ trực tuyến",
                "Periodic": null,
                "Tandc_Link": null,
                "Status": "0",
                "Introduce_product_name_en": "Meet customers' demand for money using periodic interest",
                "Introduce_product_name_vn": "Đáp ứng nhu cầu sử dụng tiền của Khách Hàng theo định kỳ lĩnh lãi",
                "Link_detail_product": "https://www.vpbank.com.vn/san-pham/gui-tiet-kiem/tiet-kiem-linh-lai-dinh-ky"
            },
            {
                "__metadata": {
                    "id": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(4)",
                    "uri": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(4)",
                    "type": "com.sap.banking.custom.onlinesaving.endpoint.v1_0.beans.OnlineSavingAccountType"
                },
                "Account_Type_ID": 4,
                "Name_EN": "Tiết kiệm easy saving trực tuyến",
                "Name_VN": "Tiết kiệm easy saving trực tuyến",
                "Periodic": null,
                "Tandc_Link": null,
                "Status": "0",
                "Introduce_product_name_en": "Flexible and profitable submissions for small funds",
                "Introduce_product_name_vn": "Gửi góp linh hoạt và sinh lời hiệu quả cho khoản tiền nhỏ nhàn rỗi.",
                "Link_detail_product": "https://www.vpbank.com.vn/san-pham/gui-tiet-kiem/tiet-kiem-gui-gop-linh-hoat-easy-savings"
            },
            {
                "__metadata": {
                    "id": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(5)",
                    "uri": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(5)",
                    "type": "com.sap.banking.custom.onlinesaving.endpoint.v1_0.beans.OnlineSavingAccountType"
                },
                "Account_Type_ID": 5,
                "Name_EN": "Tiết kiệm trả lãi trước trực tuyến",
                "Name_VN": "Tiết kiệm trả lãi trước trực tuyến",
                "Periodic": null,
                "Tandc_Link": null,
                "Status": "0",
                "Introduce_product_name_en": "Customer will receive prepaid interest amount as soon as open saving",
                "Introduce_product_name_vn": "Khách Hàng được nhận tiền lãi ngay khi gửi tiết kiệm",
                "Link_detail_product": "https://www.vpbank.com.vn/san-pham/gui-tiet-kiem/tiet-kiem-tra-lai-truoc"
            },
            {
                "__metadata": {
                    "id": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(6)",
                    "uri": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(6)",
                    "type": "com.sap.banking.custom.onlinesaving.endpoint.v1_0.beans.OnlineSavingAccountType"
                },
                "Account_Type_ID": 6,
                "Name_EN": "Tiết kiệm Phát Lộc Thịnh Vượng trực tuyến",
                "Name_VN": "Tiết kiệm Phát Lộc Thịnh Vượng trực tuyến",
                "Periodic": null,
                "Tandc_Link": "http://www.vpbank.com.vn/",
                "Status": "0",
                "Introduce_product_name_en": null,
                "Introduce_product_name_vn": null,
                "Link_detail_product": null
            },
            {
                "__metadata": {
                    "id": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(7)",
                    "uri": "https://smp-srv:8081/cb/odata/services/onlinesavingservice/OnlineSavingAccountTypes(7)",
                    "type": "com.sap.banking.custom.onlinesaving.endpoint.v1_0.beans.OnlineSavingAccountType"
                },
                "Account_Type_ID": 7,
                "Name_EN": "Tiết kiệm bảo chứng thấu chi",
                "Name_VN": "Tiết kiệm bảo chứng thấu chi",
                "Periodic": null,
                "Tandc_Link": "http://www.vpbank.com.vn/sites/default/files/pictures/Ti%E1%BA%BFt%20ki%E1%BB%87m%20th%C6%B0%E1%BB%9Dng.pdf",
                "Status": "0",
                "Introduce_product_name_en": "Combinations packgage of Online normal saving and secured overdraft product.",
                "Introduce_product_name_vn": "Khách Hàng không cần tất toán trước hạn vẫn được sử dụng tiền ngay khi cần",
                "Link_detail_product": "https://www.vpbank.com.vn/san-pham/gui-tiet-kiem/tiet-kiem-bao-chung-thau-chi"
            }
        ]
    }
}
"@
$ws.Range("G15").Value = $g15Text

# I15: mark as Done (same as H15)
$ws.Range("I15").Value = "Done"

# Update the visible selection to G15 (matches the saved cursor position)
$ws.Range("G15").Select()
